# Atualização automática da planilha
# Adds a new team member ("Tadashi Sato") to the "Organograma" sheet, as a
# new row 34 (Time Projeto / Key User N2 / Tadashi Sato / SLO - Terceiro /
# Jurídico), pushing the former rows 34-42 down to 35-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

# Remember the sheet that is active before we touch anything, so the
# workbook's active-tab selection is left exactly as it was.
$origSheet = $wb.ActiveSheet

# Insert a new blank row at position 34 - this pushes rows 34:42 down to
# 35:43 and keeps every other row (and every other sheet) untouched.
$ws.Rows.Item(34).Insert()

# The new row should look like the row directly above it (row 33, which is
# the other "Key User N2 / Jurídico" entry), so copy its formatting down.
$ws.Range("A33:E33").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

# Fill in the new team member's data.
$ws.Cells.Item(34, 1).Value = "Time Projeto"
$ws.Cells.Item(34, 2).Value = "Key User N2"
$ws.Cells.Item(34, 3).Value = "Tadashi Sato"
$ws.Cells.Item(34, 4).Value = "SLO - Terceiro"
$ws.Cells.Item(34, 5).Value = "Jurídico"

# Match the author's final cursor/scroll position on the Organograma sheet
# (frozen-pane view scrolled to row 33, cell F34 selected) ...
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$ws.Range("F34").Select()

# ... then restore whichever sheet was active before, so the workbook as a
# whole still opens on the same tab it did before this edit.
$origSheet.Activate()
